$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (stateness), shifting stateness/
# successful_transition/approach one column to the right, and fill it
# with the "proxy_score" weighted-average values.
$ws.Range("C1").EntireColumn.Insert()

# Copy the header style from the old header cell (now D1) to the new C1
$ws.Range("D1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C1").Value = "proxy_score"

$proxyScores = @(
    7.661515663058537,
    4.101782021109108,
    7.414594319451427,
    3.886213278251039,
    4.829586975558342,
    3.336017998075534,
    6.690406117635277
)

$stateness = @(
    54.72511187898955,
    29.2984430079222,
    52.96138799608162,
    27.75866627322171,
    34.49704982541673,
    23.82869998625382,
    47.78861512596627
)

for ($i = 0; $i -lt $proxyScores.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $proxyScores[$i]
    $ws.Cells.Item($row, 4).Value = $stateness[$i]
}

$ws.Range("A1").Select()
